$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.853.96'
$ws.Cells.Item(2, 5).Value = '  -1.15%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.563.89'

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'205.98"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.35%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -1.12%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.06%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'21.76"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -2.26%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.46%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -1.27%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0864"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.47%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.785.98'
$ws.Cells.Item(12, 5).Value = '  +0.05%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.562.48'
$ws.Cells.Item(13, 5).Value = '  +0.41%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'3.73"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -1.14%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -0.13%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '26.863.35'
$ws.Cells.Item(16, 5).Value = '  -1.06%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'61.29"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -2.53%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'215.22"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +0.91%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'7.36"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +2.02%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.20%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.21%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -2.22%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'2.01"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +1.73%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +1.01%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'6.74"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +2.57%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'14.88"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.02%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.08%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -1.04%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'0.0467"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +0.64%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -3.82%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.10%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.403.90'
$ws.Cells.Item(33, 5).Value = '  +1.52%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.90%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -1.15%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.45%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.918"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -2.37%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.58%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.528"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +2.20%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.810"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.60%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.07%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.71%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'5.45"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +4.41%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44, 4).Value = "'1.76"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -1.14%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'MXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(45, 4).Value = "'2.17"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.20%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'63.36"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.02%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '1.699.25'
$ws.Cells.Item(47, 5).Value = '  +0.15%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'86.58"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.13%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +2.82%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '0.0₇0978'
$ws.Cells.Item(50, 5).Value = '  -1.53%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.0948"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.58%  '
